$d = $word.ActiveDocument

# 1. Update the "Ejercicio 2" title: drop the "(Aplicación de SOLID)" suffix,
#    keep a trailing space after "Intercambiables". Only touch the run that
#    holds this text (not the preceding " 2" run) so the run split is kept
#    intact.
$d.Content.Find.Execute(
    "(Aplicación de SOLID)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2)

# 2. Simplify the closing sentence of the second exercise paragraph: remove
#    the ", aplicando los principios SOLID" clause (including the bold
#    "SOLID" run) and leave a period right after "mantener".
$d.Content.Find.Execute(
    ", aplicando los principios SOLID.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ".",
    2)

# 3. Mark the "Default Paragraph Font" style as semi-hidden, matching the
#    styles.xml change.
$style = $d.Styles("Fuentedeprrafopredeter")
$style.SemiHidden = $true
